$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "45.626.15"
$ws.Range("E2").Value = "  -2.06%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.407.57"
$ws.Range("E3").Value = "  +4.77%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.01%  "

# Row 5 - BNB
$ws.Range("D5").Value = "299.89"
$ws.Range("E5").Value = "  -1.52%  "

# Row 6 - Solana
$ws.Range("D6").Value = "97.54"
$ws.Range("E6").Value = "  -3.01%  "

# Row 7 - XRP
$ws.Range("D7").Value = "0.564"
$ws.Range("E7").Value = "  -0.39%  "

# Row 9 - Cardano
$ws.Range("D9").Value = "0.512"
$ws.Range("E9").Value = "  -2.45%  "

# Row 10 - Avalanche
$ws.Range("D10").Value = "34.92"

# Row 11 - Dogecoin
$ws.Range("D11").Value = "0.0792"
$ws.Range("E11").Value = "  +0.39%  "

# Row 12 - Polkadot
$ws.Range("E12").Value = "  -2.98%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  +0.74%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "2.761.97"
$ws.Range("E14").Value = "  +4.27%  "

# Row 15 - WrappedEther
$ws.Range("D15").Value = "2.420.17"
$ws.Range("E15").Value = "  +5.13%  "

# Row 16 - Polygon
$ws.Range("D16").Value = "0.846"
$ws.Range("E16").Value = "  +3.82%  "

# Row 17 - Chainlink
$ws.Range("D17").Value = "14.20"
$ws.Range("E17").Value = "  +2.56%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "45.598.64"
$ws.Range("E18").Value = "  -2.09%  "

# Row 19 - InternetComputer(DFINITY)
$ws.Range("D19").Value = "13.12"
$ws.Range("E19").Value = "  +0.47%  "

# Row 20 - ShibaInu
$ws.Range("D20").Value = "0.0₃0950"
$ws.Range("E20").Value = "  +1.26%  "

# Row 21 - Uniswap
$ws.Range("D21").Value = "6.23"
$ws.Range("E21").Value = "  +3.37%  "

# Row 22 - Litecoin
$ws.Range("D22").Value = "67.15"
$ws.Range("E22").Value = "  +1.40%  "

# Row 23 - BitcoinCash
$ws.Range("D23").Value = "243.06"
$ws.Range("E23").Value = "  -2.30%  "

# Row 24 - PancakeSwap
$ws.Range("D24").Value = "2.81"
$ws.Range("E24").Value = "  -3.33%  "

# Row 25 - Dai
$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  -0.21%  "

# Row 26 - ImmutableX
$ws.Range("D26").Value = "1.93"
$ws.Range("E26").Value = "  +0.26%  "

# Row 27 - InjectiveProtocol
$ws.Range("D27").Value = "38.49"
$ws.Range("E27").Value = "  -9.77%  "

# Row 28 - Toncoin
$ws.Range("D28").Value = "2.22"
$ws.Range("E28").Value = "  -2.20%  "

# Row 29 - Cosmos
$ws.Range("D29").Value = "9.81"
$ws.Range("E29").Value = "  -0.80%  "

# Row 30 - LidoDAOToken
$ws.Range("D30").Value = "3.83"
$ws.Range("E30").Value = "  +16.50%  "

# Row 31 - EthereumClassic
$ws.Range("D31").Value = "21.28"
$ws.Range("E31").Value = "  +6.34%  "

# Row 32 - WEMIXToken
$ws.Range("E32").Value = "  -1.51%  "

# Row 33 - Filecoin
$ws.Range("D33").Value = "5.55"
$ws.Range("E33").Value = "  -1.90%  "

# Row 34 - Monero
$ws.Range("D34").Value = "148.40"
$ws.Range("E34").Value = "  +0.31%  "

# Row 35 - Hedera
$ws.Range("D35").Value = "0.0776"
$ws.Range("E35").Value = "  -2.43%  "

# Row 36 - was ARBITRUM, now Kaspa
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").Value = "0.114"
$ws.Range("E36").Value = "  -0.26%  "

# Row 37 - was Kaspa, now ARBITRUM
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").Value = "1.97"
$ws.Range("E37").Value = "  +11.01%  "

# Row 38 - Stellar
$ws.Range("E38").Value = "  -1.33%  "

# Row 39 - Celestia
$ws.Range("D39").Value = "15.17"
$ws.Range("E39").Value = "  -5.39%  "

# Row 40 - RenderToken
$ws.Range("D40").Value = "3.88"
$ws.Range("E40").Value = "  -2.66%  "

# Row 41 - VeChain
$ws.Range("D41").Value = "0.0301"
$ws.Range("E41").Value = "  -0.54%  "

# Row 42 - NEARProtocol
$ws.Range("E42").Value = "  -2.36%  "

# Row 43 - Maker
$ws.Range("D43").Value = "1.953.63"
$ws.Range("E43").Value = "  +7.34%  "

# Row 44 - FirstDigitalUSD
$ws.Range("D44").Value = "0.999"
$ws.Range("E44").Value = "  +0.02%  "

# Row 45 - BitcoinSV
$ws.Range("D45").Value = "91.71"
$ws.Range("E45").Value = "  +3.98%  "

# Row 46 - Stacks
$ws.Range("D46").Value = "1.75"
$ws.Range("E46").Value = "  -11.61%  "

# Row 47 - FraxShare
$ws.Range("D47").Value = "8.70"
$ws.Range("E47").Value = "  +10.21%  "

# Row 48 - was EnergySwap, now Aave
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").Value = "102.84"
$ws.Range("E48").Value = "  +7.01%  "

# Row 49 - was Aave, now EnergySwap
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "15.31"
$ws.Range("E49").Value = "  +14.37%  "

# Row 50 - Algorand
$ws.Range("E50").Value = "  -3.53%  "

# Row 51 - RocketPoolETH
$ws.Range("D51").Value = "2.646.26"
$ws.Range("E51").Value = "  +4.79%  "
